# Updated cryptos list on Sat Jan 27 14:32:22 UTC 2024 with GitHub Actions
# Refresh per-row Price (D) / Volume(1h) (E) figures; rows 26/27 and 30/31
# also swap rank order (Dai<->ImmutableX, Toncoin<->InjectiveProtocol).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '41.777.25'
$ws.Range("E2").Value = '  +1.46%  '
# Row 3: Ethereum
$ws.Range("D3").Value = '2.270.03'
$ws.Range("E3").Value = '  +0.99%  '
# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.02%  '
# Row 5: BNB
$ws.Range("D5").Value = "'303.38"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +0.31%  '
# Row 6: Solana
$ws.Range("D6").Value = "'92.59"
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +1.55%  '
# Row 7: XRP
$ws.Range("E7").Value = '  +1.76%  '
# Row 8: USDC
$ws.Range("E8").Value = '  -0.05%  '
# Row 9: Cardano
$ws.Range("D9").Value = "'0.484"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +0.14%  '
# Row 10: Avalanche
$ws.Range("D10").Value = "'32.53"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  +1.85%  '
# Row 11: OKB
$ws.Range("D11").Value = "'53.28"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -0.75%  '
# Row 12: Dogecoin
$ws.Range("E12").Value = '  +0.41%  '
# Row 13: TRON
$ws.Range("E13").Value = '  -1.40%  '
# Row 14: Polkadot
$ws.Range("E14").Value = '  +1.48%  '
# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '2.621.63'
$ws.Range("E15").Value = '  +0.98%  '
# Row 16: Chainlink
$ws.Range("D16").Value = "'14.26"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  +1.19%  '
# Row 17: WrappedEther
$ws.Range("D17").Value = '2.273.93'
$ws.Range("E17").Value = '  +1.32%  '
# Row 18: Polygon
$ws.Range("D18").Value = "'0.776"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  +3.86%  '
# Row 19: WrappedBTC
$ws.Range("D19").Value = '41.685.15'
$ws.Range("E19").Value = '  +1.40%  '
# Row 20: InternetComputer(DFINITY)
$ws.Range("D20").Value = "'12.47"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +2.97%  '
# Row 21: ShibaInu
$ws.Range("E21").Value = '  +0.33%  '
# Row 22: Uniswap
$ws.Range("E22").Value = '  +1.43%  '
# Row 23: Litecoin
$ws.Range("D23").Value = "'67.10"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +0.70%  '
# Row 24: BitcoinCash
$ws.Range("D24").Value = "'239.96"
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -0.12%  '
# Row 25: PancakeSwap
$ws.Range("D25").Value = "'2.59"
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  +0.95%  '
# Row 26: ImmutableX
$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").Value = "'1.93"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +4.55%  '
# Row 27: Dai
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +0.02%  '
# Row 28: EthereumClassic
$ws.Range("E28").Value = '  +1.28%  '
# Row 29: Cosmos
$ws.Range("D29").Value = "'9.54"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -0.68%  '
# Row 30: InjectiveProtocol
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = "'35.69"
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  +6.98%  '
# Row 31: Toncoin
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = "'2.07"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -5.01%  '
# Row 33: Filecoin
$ws.Range("D33").Value = "'5.24"
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  +1.76%  '
# Row 34: FirstDigitalUSD
$ws.Range("E34").Value = '  -0.05%  '
# Row 35: Hedera
$ws.Range("D35").Value = "'0.0745"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  +1.71%  '
# Row 36: LidoDAOToken
$ws.Range("E36").Value = '  -0.74%  '
# Row 37: Celestia
$ws.Range("E37").Value = '  +1.43%  '
# Row 38: WEMIXToken
$ws.Range("E38").Value = '  +0.57%  '
# Row 39: Kaspa
$ws.Range("E39").Value = '  +1.95%  '
# Row 40: Stellar
$ws.Range("E40").Value = '  +1.03%  '
# Row 41: ARBITRUM
$ws.Range("E41").Value = '  +0.94%  '
# Row 42: RenderToken
$ws.Range("E42").Value = '  -0.10%  '
# Row 43: Maker
$ws.Range("D43").Value = '2.002.49'
$ws.Range("E43").Value = '  -3.01%  '
# Row 44: EnergySwap
$ws.Range("D44").Value = "'19.15"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  -5.36%  '
# Row 45: VeChain
$ws.Range("E45").Value = '  +2.03%  '
# Row 46: FraxShare
$ws.Range("D46").Value = "'10.34"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +1.28%  '
# Row 47: ApeXProtocol
$ws.Range("D47").Value = "'2.12"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +4.03%  '
# Row 48: NEARProtocol
$ws.Range("D48").Value = "'2.89"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -1.70%  '
# Row 49: MultiversX
$ws.Range("D49").Value = "'52.57"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +3.50%  '
# Row 50: TrustWalletToken
$ws.Range("E50").Value = '  +1.06%  '
# Row 51: Stacks
$ws.Range("E51").Value = '  +0.84%  '
